# Applies the "Diablando A Temporada" content/formatting edits described
# by the commit diff:
#  1. Rewrites the "Objetivo" body paragraph with new copy, tags it with
#     the Ttulo1 (Heading 1) style (visually overridden back to normal
#     run formatting + a muted grey theme colour), and splits the text
#     into two runs (body + trailing period).
#  2. Underlines the blank paragraph mark right under "Premissas e
#     Restrições".
#  3. Rewrites the "Plataforma Online Integrada" bullet with new,
#     shorter copy, and restructures the bullets that follow it.

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Objetivo" section body paragraph
# ---------------------------------------------------------------------
$objetivoPara = $d.Paragraphs(11)
if ($objetivoPara.Range.Text -notmatch "O objetivo do site") {
    throw "Paragraph 11 is not the 'Objetivo' body paragraph any more"
}

$objetivoXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Ttulo1"/>
    <w:ind w:firstLine="708"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="595959" w:themeColor="text1" w:themeTint="A6"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="595959" w:themeColor="text1" w:themeTint="A6"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>O objetivo principal do site é fornecer uma plataforma centralizada e acessível para os fãs de Diablo 3 se conectarem, trocarem conhecimentos e experiências, e se envolverem em discussões e atividades relacionadas ao jogo. O site deve ser um espaço acolhedor e informativo para jogadores de todos os níveis, desde iniciantes até jogadores avançados</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Arial" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="595959" w:themeColor="text1" w:themeTint="A6"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
'@
$objetivoPara.Range.InsertXML($objetivoXml)

# ---------------------------------------------------------------------
# 2) Blank paragraph right after "Premissas e Restrições" -> underline
#    the paragraph mark.
# ---------------------------------------------------------------------
$restricoesHeading = $d.Paragraphs(24)
if ($restricoesHeading.Range.Text -notmatch "Premissas e Restri") {
    throw "Paragraph 24 is not the 'Premissas e Restrições' heading any more"
}
$blankAfterHeading = $d.Paragraphs(25)
if ($blankAfterHeading.Range.Text.Trim().Length -ne 0) {
    throw "Paragraph 25 is not the expected blank paragraph"
}

$blankXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@
$blankAfterHeading.Range.InsertXML($blankXml)

# ---------------------------------------------------------------------
# 3) "Plataforma Online Integrada" bullet + the two bullets that follow
#    it get rewritten/restructured.
# ---------------------------------------------------------------------
$plataformaPara = $d.Paragraphs(26)
if ($plataformaPara.Range.Text -notmatch "Plataforma Online Integrada") {
    throw "Paragraph 26 is not the 'Plataforma Online Integrada' bullet any more"
}

# Replace that single bullet with:
#   - the rewritten "Conteúdo abrangente e atualizado" bullet
#   - a new blank (non-numbered) PargrafodaLista paragraph
$plataformaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="15"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Conteúdo abrangente e atualizado: O site deve oferecer um conteúdo rico e atualizado sobre Diablo 3, incluindo guias, dicas, notícias, atualizações do jogo, informações sobre personagens, builds, estratégias de jogo e outros recursos relevantes. O objetivo é fornecer aos usuários acesso a informações completas e confiáveis para melhorar sua experiência de jogo.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@
$plataformaPara.Range.InsertXML($plataformaXml)

# The old blank bullet-style paragraph that used to sit between the
# "Plataforma..." bullet and the "Conteúdo e Atualizações..." bullet is
# now redundant (we already inserted its replacement above), so delete
# it: after the insert above, paragraph 26 is the rewritten bullet,
# paragraph 27 is the freshly-inserted blank paragraph, and paragraph 28
# is the old blank paragraph that needs to go away.
$oldBlank = $d.Paragraphs(28)
if ($oldBlank.Range.Text.Trim().Length -ne 0) {
    throw "Paragraph 28 is not the stale blank paragraph"
}
$oldBlank.Range.Delete() | Out-Null

# ---------------------------------------------------------------------
# "Conteúdo e Atualizações Relevantes" bullet gets new copy (split into
# two runs: body + "s."), followed by a new blank PargrafodaLista
# paragraph, a new numbered "Ambiente inclusivo e respeitoso" bullet and
# a new blank numbered bullet.
# ---------------------------------------------------------------------
$conteudoPara = $d.Paragraphs(28)
if ($conteudoPara.Range.Text -notmatch "Conte.do e Atualiza") {
    throw "Paragraph 28 is not the 'Conteúdo e Atualizações Relevantes' bullet any more"
}

$conteudoXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="15"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Comunicação e interação: O site deve promover a interação entre os membros da comunidade, oferecendo recursos como fóruns de discussão, chat em tempo real, grupos de jogadores, sistema de mensagens privadas e outras ferramentas de comunicação. A interação entre os jogadores é fundamental para criar um senso de comunidade e compartilhamento de conhecimento</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>s.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="15"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Ambiente inclusivo e respeitoso: Uma das premissas mais importantes é criar um ambiente inclusivo, respeitoso e livre de discriminação, onde todos os membros da comunidade se sintam bem-vindos e confortáveis para participar. Regras claras de conduta devem ser estabelecidas e medidas devem ser tomadas para evitar comportamentos inadequados, assédio ou qualquer forma de discriminação. O objetivo é criar uma comunidade saudável e amigável, onde os jogadores possam interagir de maneira positiva.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="15"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
</w:p>
'@
$conteudoPara.Range.InsertXML($conteudoXml)

Write-Host "Edits applied."
